$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @(newPrice, newVolume)  ($null means "leave unchanged")
$updates = @{
    2  = @("30.391.44", "  +1.24%  ")
    3  = @("1.922.45", "  +0.66%  ")
    4  = @("1.001", "  +0.06%  ")
    5  = @("0.8094", "  +2.68%  ")
    6  = @("244.82", "  +1.07%  ")
    7  = @("1.000", "  +0.04%  ")
    8  = @("0.3254", "  +2.64%  ")
    9  = @("27.35", "  +3.56%  ")
    10 = @("0.07294", "  +6.00%  ")
    11 = @("0.7927", "  +6.57%  ")
    12 = @("0.08101", "  +1.42%  ")
    13 = @("1.921.25", "  +0.64%  ")
    14 = @("5.427", "  +4.32%  ")
    15 = @("94.74", "  +1.85%  ")
    16 = @("30.379.80", "  +1.23%  ")
    17 = @("14.38", "  +3.18%  ")
    18 = @("6.084", "  +3.54%  ")
    19 = @("254.22", "  +3.46%  ")
    20 = @("0.000007865", "  +1.48%  ")
    21 = @("2.177.64", "  +0.81%  ")
    23 = @("8.046", "  +17.46%  ")
    24 = @($null, "  +0.09%  ")
    25 = @("0.1641", "  +18.38%  ")
    26 = @("9.558", "  +3.44%  ")
    27 = @("167.50", "  -0.45%  ")
    28 = @("19.17", "  +1.53%  ")
    29 = @("2.151", "  +5.49%  ")
    30 = @("1.377", "  +0.67%  ")
    31 = @($null, "  +1.31%  ")
    32 = @("4.359", "  +0.90%  ")
    33 = @("4.157", "  +1.68%  ")
    34 = @("0.05650", "  +1.84%  ")
    35 = @("1.305", "  +3.99%  ")
    36 = @("0.7454", "  +1.55%  ")
    37 = @("1.003", "  +0.45%  ")
    38 = @("2.721", "  +0.02%  ")
    39 = @("0.01960", "  +1.53%  ")
    40 = @("2.813", "  +0.94%  ")
    41 = @("0.4507", "  +1.89%  ")
    42 = @("73.86", "  +2.19%  ")
    43 = @("5.998", "  -2.47%  ")
    44 = @("1.942", "  +3.46%  ")
    45 = @("0.8550", "  +1.78%  ")
    46 = @($null, "  +0.05%  ")
    47 = @("103.53", "  +3.04%  ")
    48 = @("1.032.01", "  +4.79%  ")
    49 = @("9.989", "  +2.67%  ")
    50 = @("7.669", "  +1.73%  ")
    51 = @("2.073.68", "  +1.06%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $price = $vals[0]
    $volume = $vals[1]

    if ($null -ne $price) {
        $cell = $ws.Cells.Item($row, 4)
        # Force the price to be stored as text (it looks numeric, e.g. "1.001"),
        # matching the original inline-string cell content, then restore the
        # cell's style so no stray formatting is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $price
        $cell.Style = "Normal"
    }
    if ($null -ne $volume) {
        $ws.Cells.Item($row, 5).Value = $volume
    }
}
